# Assignment2/Benchmarks.xlsx - "Change probability assignment approach and rerun benchmark tests"
#
# The underlying test data was re-generated (new probability assignment
# approach for distributing operations across threads) and the benchmark
# was rerun on different hardware, so:
#   - the header labels for the three tests were reworded / reordered
#   - the hardware note text was updated
#   - all of the recorded timing numbers (columns B:G, rows 3-9) changed
#   - the active selection moved to B14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels (row 1, merged cells) -----------------------------
# Order the writes so the shared-string table is rebuilt in the same
# sequence the author's copy of Excel produced it in.
$ws.Range("B1").Value = "Test 1 (mixed)"
$ws.Range("D1").Value = "Test 2 (write-dominated)"

# --- Notes section (rows 11-13) ---------------------------------------
$ws.Range("B13").Value = "Computed using 3.1 GHz Intel Core i7 quad-core processor (4 physical, 8 logical)"

$ws.Range("F1").Value = "Test 3 (read-dominated)"

# --- Benchmark results grid (A3:G9) -----------------------------------
# Row for concurrency level 1
$ws.Range("B3").Value = 48193
$ws.Range("C3").Value = 31500
$ws.Range("D3").Value = 77836
$ws.Range("E3").Value = 49608
$ws.Range("F3").Value = 18089
$ws.Range("G3").Value = 10130

# Row for concurrency level 2
$ws.Range("B4").Value = 24859
$ws.Range("C4").Value = 33891
$ws.Range("D4").Value = 40038
$ws.Range("E4").Value = 53456
$ws.Range("F4").Value = 9499
$ws.Range("G4").Value = 12097

# Row for concurrency level 4
$ws.Range("B5").Value = 13616
$ws.Range("C5").Value = 33402
$ws.Range("D5").Value = 21839
$ws.Range("E5").Value = 53193
$ws.Range("F5").Value = 5229
$ws.Range("G5").Value = 11814

# Row for concurrency level 8
$ws.Range("B6").Value = 8894
$ws.Range("C6").Value = 33360
$ws.Range("D6").Value = 13515
$ws.Range("E6").Value = 53357
$ws.Range("F6").Value = 3931
$ws.Range("G6").Value = 11978

# Row for concurrency level 16
$ws.Range("B7").Value = 8809
$ws.Range("C7").Value = 34403
$ws.Range("D7").Value = 13439
$ws.Range("E7").Value = 53284
$ws.Range("F7").Value = 3881
$ws.Range("G7").Value = 12494

# Row for concurrency level 32
$ws.Range("B8").Value = 8627
$ws.Range("C8").Value = 34103
$ws.Range("D8").Value = 13227
$ws.Range("E8").Value = 53641
$ws.Range("F8").Value = 4001
$ws.Range("G8").Value = 12622

# Row for concurrency level 64
$ws.Range("B9").Value = 8830
$ws.Range("C9").Value = 34062
$ws.Range("D9").Value = 13200
$ws.Range("E9").Value = 53659
$ws.Range("F9").Value = 3956
$ws.Range("G9").Value = 12714

# --- Selection moved to B14 -------------------------------------------
[void]$ws.Range("B14").Select()
